$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the join conditions in every SQL query cell: the queries used to join
# on the generic "id" column; they now join on the fully-qualified
# "study_id" / "participant_id" columns.
# Cell processing order matches the order the original authors' save produced
# in the shared-string table (Count query / Study query / Participant query /
# Diagnosis query / Treatment query / TreatmentResponse query / Survival query).
$cellsToFix = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cellsToFix) {
    $sql = $ws.Range($addr).Value2
    $sql = $sql.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $sql = $sql.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $sql = $sql.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $sql = $sql.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $sql = $sql.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $sql = $sql.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    $ws.Range($addr).Value = $sql
}

# The "StatQuery" cell (C2) now uses a smaller, normal-weight 11pt font
# instead of the 12pt font used elsewhere, so give it its own font.
$ws.Range("C2").Font.Size = 11
$ws.Range("C2").Font.ThemeColor = 1

# Column C was widened (and is no longer an auto "best fit" column).
$ws.Columns.Item(3).ColumnWidth = 66
